$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.063.64'
$ws.Range('E2').Value = '  -3.35%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.536.30'
$ws.Range('E3').Value = '  -3.79%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.23'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.53'
$ws.Range('E6').Value = '  -3.24%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.532.81'
$ws.Range('E7').Value = '  -3.81%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('E9').Value = '  -2.44%  '

$ws.Range('E10').Value = '  -2.09%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.84'
$ws.Range('E11').Value = '  -3.25%  '

$ws.Range('E12').Value = '  -3.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('E13').Value = '  -4.46%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.130.04'
$ws.Range('E14').Value = '  -3.80%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.96'
$ws.Range('E15').Value = '  -2.10%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.529.34'
$ws.Range('E16').Value = '  -4.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.053.71'
$ws.Range('E17').Value = '  -3.34%  '

$ws.Range('E18').Value = '  +0.90%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.38'
$ws.Range('E19').Value = '  -1.91%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.45'
$ws.Range('E20').Value = '  -3.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '451.85'
$ws.Range('E21').Value = '  -2.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.38'
$ws.Range('E22').Value = '  -5.28%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.639'
$ws.Range('E23').Value = '  -1.40%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.14'
$ws.Range('E24').Value = '  -0.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.672.88'
$ws.Range('E25').Value = '  -3.85%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000124'
$ws.Range('E27').Value = '  -1.52%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.27'
$ws.Range('E28').Value = '  -5.69%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.32'
$ws.Range('E29').Value = '  -8.09%  '

$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.69'
$ws.Range('E30').Value = '  -0.83%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.55'
$ws.Range('E31').Value = '  -2.82%  '

$ws.Range('E32').Value = '  -0.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.96'
$ws.Range('E33').Value = '  -3.41%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.90'
$ws.Range('E34').Value = '  -5.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.19'
$ws.Range('E35').Value = '  -4.01%  '

$ws.Range('E36').Value = '  -5.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.528.99'
$ws.Range('E37').Value = '  -3.81%  '

$ws.Range('E38').Value = '  -4.20%  '

$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  -0.21%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '176.07'
$ws.Range('E41').Value = '  -0.96%  '

$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.15'
$ws.Range('E42').Value = '  -1.85%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.61'
$ws.Range('E43').Value = '  -4.90%  '

$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0877'
$ws.Range('E44').Value = '  -2.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.892'
$ws.Range('E45').Value = '  -3.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.79'
$ws.Range('E46').Value = '  -2.03%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.41'
$ws.Range('E47').Value = '  +3.94%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.71'
$ws.Range('E48').Value = '  -0.71%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.23'
$ws.Range('E49').Value = '  -1.29%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.65'
$ws.Range('E50').Value = '  -2.34%  '

$ws.Range('E51').Value = '  -3.28%  '
